$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Definition-Name + Description updated
$ws.Range("A2").Value = "D0.18565903795358007"
$ws.Range("D2").Value = "Mon, 19 Dec 2022 05:44:49 -0800"

# Row 3: Definition-Name + Description updated
$ws.Range("A3").Value = "D0.3123609309210864"
$ws.Range("D3").Value = "Mon, 19 Dec 2022 05:44:49 -0800"

# Row 4: Definition-Name + Description updated
$ws.Range("A4").Value = "D0.6471154244116394"
$ws.Range("D4").Value = "Mon, 19 Dec 2022 05:44:49 -0800"
